$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.933.50"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "2.618.74"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'529.34"
$ws.Range("E5").Value = "  +3.80%  "
$ws.Range("D6").Value = "'154.74"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.592"
$ws.Range("E8").Value = "  +1.37%  "
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").Value = "'0.109"
$ws.Range("E10").Value = "  +4.90%  "
$ws.Range("D11").Value = "'0.348"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "3.079.13"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").Value = "60.938.68"
$ws.Range("E14").Value = "  +0.74%  "
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "'0.0000144"
$ws.Range("E16").Value = "  +2.65%  "
$ws.Range("D17").Value = "2.625.06"
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").Value = "'4.78"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").Value = "'354.05"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("D20").Value = "'10.61"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("E21").Value = "  +1.79%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "'61.50"
$ws.Range("E23").Value = "  +1.60%  "
$ws.Range("D24").Value = "'0.430"
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("D25").Value = "'0.169"
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("D26").Value = "2.735.59"
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.34%  "
$ws.Range("D28").Value = "0.0₃0866"
$ws.Range("E28").Value = "  +2.27%  "
$ws.Range("D29").Value = "'7.38"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  +7.47%  "
$ws.Range("D32").Value = "'19.46"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  +2.76%  "
$ws.Range("D34").Value = "'150.89"
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("D35").Value = "'4.17"
$ws.Range("E35").Value = "  +3.01%  "
$ws.Range("E36").Value = "  +1.01%  "
$ws.Range("D37").Value = "'0.929"
$ws.Range("E37").Value = "  +9.54%  "
$ws.Range("D38").Value = "'0.892"
$ws.Range("E38").Value = "  +2.79%  "
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("D41").Value = "'295.16"
$ws.Range("E41").Value = "  -1.72%  "
$ws.Range("D42").Value = "'0.638"
$ws.Range("E42").Value = "  +2.62%  "
$ws.Range("E43").Value = "  +2.37%  "
$ws.Range("D44").Value = "'0.0561"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "'19.73"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("E47").Value = "  +2.45%  "
$ws.Range("D48").Value = "'4.90"
$ws.Range("E48").Value = "  +1.35%  "
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").Value = "'19.05"
$ws.Range("E50").Value = "  +5.78%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "'1.83"
$ws.Range("E51").Value = "  +2.52%  "
